# Arsenal_stats.xlsx update
#  1) Rename the per-category stats sheets to human-friendly, spaced names.
#  2) Bump every player's "Age" column (format "YY-DDD", years-days since
#     birthday) forward by a single day across all nine stats sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet renames
# ---------------------------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"    = "Shooting Stats"
    "PassingStats"     = "Passing Stats"
    "PassTypes"        = "Pass Types"
    "GoalShotCreation" = "Goal & Shot Creation"
    "DefensiveActions" = "Defensive Actions"
    "PlayingTime"      = "Playing Time"
    "MiscStats"        = "Miscellaneous Stats"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# ---------------------------------------------------------------------
# 2) Age column (+1 day), applied identically on every stats sheet
# ---------------------------------------------------------------------
function Get-NextAge {
    param([string]$age)

    $parts = $age -split "-"
    $years = [int]$parts[0]
    $days = [int]$parts[1]

    if ($days -eq 364) {
        $years = $years + 1
        $days = 0
    } else {
        $days = $days + 1
    }

    $daysText = $days.ToString("D3")
    return "$years-$daysText"
}

$statsSheetNames = @(
    "Standard Stats",
    "Shooting Stats",
    "Passing Stats",
    "Pass Types",
    "Goal & Shot Creation",
    "Defensive Actions",
    "Possession",
    "Playing Time",
    "Miscellaneous Stats"
)

foreach ($sheetName in $statsSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $cur = $cell.Value2
        if ($cur -match "^\d+-\d{3}$") {
            $cell.Value2 = Get-NextAge $cur
        }
    }
}
